# This script applies the numeric updates recorded in the commit diff
# ("chore: update Sheets via scheduled runner") to the Leve-profit tables
# on each class sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
#
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 68965.875
$ws.Range("I11").Value = 68965.875
$ws.Range("K11").Value = 68965.875
$ws.Range("M11").Value = -68825.875
$ws.Range("H33").Value = 34337.09
$ws.Range("I33").Value = 47057.312
$ws.Range("K33").Value = 47057.312
$ws.Range("M33").Value = -46828.312
$ws.Range("H43").Value = 1113.5714
$ws.Range("I43").Value = 1019
$ws.Range("K43").Value = 1019
$ws.Range("M43").Value = -950
$ws.Range("H86").Value = 8343723.5
$ws.Range("I86").Value = 6971.2
$ws.Range("J86").Value = 22238312
$ws.Range("K86").Value = 6971.2
$ws.Range("L86").Value = 22238312
$ws.Range("M86").Value = -5848.2
$ws.Range("N86").Value = -22240558
$ws.Range("H88").Value = 3232.3333
$ws.Range("J88").Value = 3849.5
$ws.Range("L88").Value = 3849.5
$ws.Range("N88").Value = -4661.5
$ws.Range("H89").Value = 8343723.5
$ws.Range("I89").Value = 6971.2
$ws.Range("J89").Value = 22238312
$ws.Range("K89").Value = 34856
$ws.Range("L89").Value = 111191560
$ws.Range("M89").Value = -29240
$ws.Range("N89").Value = -111202792
$ws.Range("H91").Value = 3232.3333
$ws.Range("J91").Value = 3849.5
$ws.Range("L91").Value = 3849.5
$ws.Range("N91").Value = -6657.5
$ws.Range("H106").Value = 59362.055
$ws.Range("I106").Value = 3251.2144
$ws.Range("K106").Value = 3251.2144
$ws.Range("M106").Value = -2620.2144
$ws.Range("H113").Value = 83337570
$ws.Range("I113").Value = 200002960
$ws.Range("J113").Value = 5143
$ws.Range("K113").Value = 200002960
$ws.Range("L113").Value = 5143
$ws.Range("M113").Value = -199999706
$ws.Range("N113").Value = -11651
$ws.Range("H116").Value = 53133316
$ws.Range("I116").Value = 83668170
$ws.Range("K116").Value = 83668170
$ws.Range("M116").Value = -83664728
$ws.Range("H131").Value = 13033.333
$ws.Range("I131").Value = 5000
$ws.Range("K131").Value = 15000
$ws.Range("M131").Value = -9960
$ws.Range("H137").Value = 1727.5883
$ws.Range("I137").Value = 1624.6
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 4873.799999999999
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -2323.799999999999
$ws.Range("N137").Value = -12600
$ws.Range("H138").Value = 2419.8594
$ws.Range("J138").Value = 3303.525
$ws.Range("L138").Value = 9910.575000000001
$ws.Range("N138").Value = -20190.575
$ws.Range("H141").Value = 1709.2142
$ws.Range("I141").Value = 1802.4166
$ws.Range("K141").Value = 5407.2498
$ws.Range("M141").Value = -227.2497999999996

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 1442.5
$ws.Range("I41").Value = 1442.5
$ws.Range("K41").Value = 1442.5
$ws.Range("M41").Value = -1028.5
$ws.Range("H122").Value = 12348396
$ws.Range("I122").Value = 16668937
$ws.Range("J122").Value = 3992.7144
$ws.Range("K122").Value = 50006811
$ws.Range("L122").Value = 11978.1432
$ws.Range("M122").Value = -50004361
$ws.Range("N122").Value = -16878.1432

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1461.3334
$ws.Range("I86").Value = 1176.6
$ws.Range("J86").Value = 2173.1667
$ws.Range("K86").Value = 1176.6
$ws.Range("L86").Value = 2173.1667
$ws.Range("M86").Value = -53.59999999999991
$ws.Range("N86").Value = -4419.1667
$ws.Range("H89").Value = 1461.3334
$ws.Range("I89").Value = 1176.6
$ws.Range("J89").Value = 2173.1667
$ws.Range("K89").Value = 5883
$ws.Range("L89").Value = 10865.8335
$ws.Range("M89").Value = -267
$ws.Range("N89").Value = -22097.8335
$ws.Range("H94").Value = 3103.1904
$ws.Range("J94").Value = 2610
$ws.Range("L94").Value = 2610
$ws.Range("N94").Value = -3512
$ws.Range("H105").Value = 1892.0625
$ws.Range("I105").Value = 1848.0714
$ws.Range("K105").Value = 1848.0714
$ws.Range("M105").Value = -101.0714
$ws.Range("H134").Value = 888
$ws.Range("I134").Value = 686.8570999999999
$ws.Range("K134").Value = 2060.5713
$ws.Range("M134").Value = 474.4287000000004

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2868.8196
$ws.Range("I31").Value = 966.73334
$ws.Range("K31").Value = 966.73334
$ws.Range("M31").Value = -671.73334
$ws.Range("H34").Value = 2868.8196
$ws.Range("I34").Value = 966.73334
$ws.Range("K34").Value = 966.73334
$ws.Range("M34").Value = -764.73334
$ws.Range("H58").Value = 1116.6818
$ws.Range("I58").Value = 884.1905
$ws.Range("K58").Value = 884.1905
$ws.Range("M58").Value = -681.1905
$ws.Range("H122").Value = 2028
$ws.Range("I122").Value = 2084.1
$ws.Range("K122").Value = 6252.299999999999
$ws.Range("M122").Value = -3802.299999999999
$ws.Range("H136").Value = 1116.6818
$ws.Range("I136").Value = 884.1905
$ws.Range("K136").Value = 2652.5715
$ws.Range("M136").Value = -102.5715

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1198.8235
$ws.Range("I113").Value = 579.5
$ws.Range("J113").Value = 1536.6364
$ws.Range("K113").Value = 1738.5
$ws.Range("L113").Value = 4609.9092
$ws.Range("M113").Value = 431.5
$ws.Range("N113").Value = -8949.9092
$ws.Range("H115").Value = 5000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 5000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 15000
$ws.Range("N115").Value = -17350
$ws.Range("M115").ClearContents()
$ws.Range("H118").Value = 3029.8
$ws.Range("J118").Value = 4750
$ws.Range("L118").Value = 14250
$ws.Range("N118").Value = -16736
$ws.Range("H122").Value = 840.2
$ws.Range("J122").Value = 898.5714
$ws.Range("L122").Value = 8087.1426
$ws.Range("N122").Value = -12987.1426

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 700.7692
$ws.Range("I97").Value = 668.0909
$ws.Range("J97").Value = 880.5
$ws.Range("K97").Value = 668.0909
$ws.Range("L97").Value = 880.5
$ws.Range("M97").Value = -172.0909
$ws.Range("N97").Value = -1872.5
$ws.Range("H102").Value = 1535.3478
$ws.Range("I102").Value = 1236.5555
$ws.Range("J102").Value = 2611
$ws.Range("K102").Value = 1236.5555
$ws.Range("L102").Value = 2611
$ws.Range("M102").Value = 385.4445000000001
$ws.Range("N102").Value = -5855
$ws.Range("H122").Value = 20002314
$ws.Range("I122").Value = 2283.0715
$ws.Range("K122").Value = 6849.2145
$ws.Range("M122").Value = -4399.2145
$ws.Range("H132").Value = 5171.174
$ws.Range("J132").Value = 5102
$ws.Range("L132").Value = 15306
$ws.Range("N132").Value = -20366

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2475.5557
$ws.Range("I40").Value = 2463.8333
$ws.Range("K40").Value = 2463.8333
$ws.Range("M40").Value = -2327.8333
$ws.Range("H55").Value = 308.08334
$ws.Range("I55").Value = 257
$ws.Range("J55").Value = 461.33334
$ws.Range("K55").Value = 257
$ws.Range("L55").Value = 461.33334
$ws.Range("M55").Value = -84
$ws.Range("N55").Value = -807.33334
$ws.Range("H68").Value = 21789.8
$ws.Range("I68").Value = 7983.6665
$ws.Range("J68").Value = 42499
$ws.Range("K68").Value = 7983.6665
$ws.Range("L68").Value = 42499
$ws.Range("M68").Value = -7234.6665
$ws.Range("N68").Value = -43997
$ws.Range("H70").Value = 44081.5
$ws.Range("J70").Value = 44081.5
$ws.Range("L70").Value = 44081.5
$ws.Range("N70").Value = -44621.5
$ws.Range("H71").Value = 21789.8
$ws.Range("I71").Value = 7983.6665
$ws.Range("J71").Value = 42499
$ws.Range("K71").Value = 39918.3325
$ws.Range("L71").Value = 212495
$ws.Range("M71").Value = -36174.3325
$ws.Range("N71").Value = -219983
$ws.Range("H73").Value = 44081.5
$ws.Range("J73").Value = 44081.5
$ws.Range("L73").Value = 44081.5
$ws.Range("N73").Value = -45953.5
$ws.Range("H122").Value = 3717.9546
$ws.Range("I122").Value = 2677.077
$ws.Range("J122").Value = 5221.4443
$ws.Range("K122").Value = 8031.231000000001
$ws.Range("L122").Value = 15664.3329
$ws.Range("M122").Value = -5581.231000000001
$ws.Range("N122").Value = -20564.3329
$ws.Range("H132").Value = 3361.7144
$ws.Range("I132").Value = 2249.3845
$ws.Range("J132").Value = 4619.1304
$ws.Range("K132").Value = 6748.1535
$ws.Range("L132").Value = 13857.3912
$ws.Range("M132").Value = -4218.1535
$ws.Range("N132").Value = -18917.3912

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 16917.334
$ws.Range("J41").Value = 16917.334
$ws.Range("L41").Value = 16917.334
$ws.Range("N41").Value = -17697.334
$ws.Range("H81").Value = 6670603
$ws.Range("I81").Value = 1069.4
$ws.Range("J81").Value = 20009670
$ws.Range("K81").Value = 2138.8
$ws.Range("L81").Value = 40019340
$ws.Range("M81").Value = -1077.8
$ws.Range("N81").Value = -40021462
$ws.Range("H84").Value = 6670603
$ws.Range("I84").Value = 1069.4
$ws.Range("J84").Value = 20009670
$ws.Range("K84").Value = 10694
$ws.Range("L84").Value = 200096700
$ws.Range("M84").Value = -5390
$ws.Range("N84").Value = -200107308
$ws.Range("H107").Value = 403.46155
$ws.Range("I107").Value = 386.18182
$ws.Range("K107").Value = 1158.54546
$ws.Range("M107").Value = 761.45454
$ws.Range("H132").Value = 4681.6333
$ws.Range("I132").Value = 4297.8
$ws.Range("K132").Value = 12893.4
$ws.Range("M132").Value = -10363.4
$ws.Range("H141").Value = 40000
$ws.Range("J141").Value = 40000
$ws.Range("L141").Value = 40000
$ws.Range("N141").Value = -50360
